$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.067515367707383
$ws.Cells.Item(2, 4).Value = 1.066648265962145
$ws.Cells.Item(2, 5).Value = 1.071631246931135
$ws.Cells.Item(2, 6).Value = 1.081017199031097
$ws.Cells.Item(2, 9).Value = 1.056387605843959
$ws.Cells.Item(2, 10).Value = 1.07245987330346
$ws.Cells.Item(2, 11).Value = 1.069358249899161
$ws.Cells.Item(2, 12).Value = 1.074327906346831
$ws.Cells.Item(2, 13).Value = 1.083689120700721

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.068928090650741
$ws.Cells.Item(3, 4).Value = 1.067776926209492
$ws.Cells.Item(3, 5).Value = 1.072950565568603
$ws.Cells.Item(3, 6).Value = 1.082508504673063
$ws.Cells.Item(3, 9).Value = 1.056912375903674
$ws.Cells.Item(3, 10).Value = 1.073526819397443
$ws.Cells.Item(3, 11).Value = 1.070301762635522
$ws.Cells.Item(3, 12).Value = 1.075462574503107
$ws.Cells.Item(3, 13).Value = 1.084997162449741

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.069840582125224
$ws.Cells.Item(4, 4).Value = 1.068505779158424
$ws.Cells.Item(4, 5).Value = 1.073802457860656
$ws.Cells.Item(4, 6).Value = 1.0834720702397
$ws.Cells.Item(4, 9).Value = 1.057249789916709
$ws.Cells.Item(4, 10).Value = 1.074215097625616
$ws.Cells.Item(4, 11).Value = 1.070910227125935
$ws.Cells.Item(4, 12).Value = 1.076194438545746
$ws.Cells.Item(4, 13).Value = 1.0858416276603

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.070223809267358
$ws.Cells.Item(5, 4).Value = 1.068811843006533
$ws.Cells.Item(5, 5).Value = 1.074160170210608
$ws.Cells.Item(5, 6).Value = 1.083876823210908
$ws.Cells.Item(5, 9).Value = 1.057391127720014
$ws.Cells.Item(5, 10).Value = 1.074503950797301
$ws.Cells.Item(5, 11).Value = 1.071165539234297
$ws.Cells.Item(5, 12).Value = 1.076501559590204
$ws.Cells.Item(5, 13).Value = 1.086196186617422

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.070288132461975
$ws.Cells.Item(6, 4).Value = 1.06886321233133
$ws.Cells.Item(6, 5).Value = 1.074220207063781
$ws.Cells.Item(6, 6).Value = 1.083944763876281
$ws.Cells.Item(6, 9).Value = 1.05741482907846
$ws.Cells.Item(6, 10).Value = 1.074552421413199
$ws.Cells.Item(6, 11).Value = 1.071208378877077
$ws.Cells.Item(6, 12).Value = 1.076553094219964
$ws.Cells.Item(6, 13).Value = 1.086255692184815

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.069845704326204
$ws.Cells.Item(7, 4).Value = 1.068509870151679
$ws.Cells.Item(7, 5).Value = 1.073807239285193
$ws.Cells.Item(7, 6).Value = 1.083477479856299
$ws.Cells.Item(7, 9).Value = 1.057251680483692
$ws.Cells.Item(7, 10).Value = 1.07421895924976
$ws.Cells.Item(7, 11).Value = 1.070913640524048
$ws.Cells.Item(7, 12).Value = 1.076198544483796
$ws.Cells.Item(7, 13).Value = 1.085846367070904

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.067993146799554
$ws.Cells.Item(8, 4).Value = 1.067030008455863
$ws.Cells.Item(8, 5).Value = 1.072077492496161
$ws.Cells.Item(8, 6).Value = 1.081521488884435
$ws.Cells.Item(8, 9).Value = 1.056565400780958
$ws.Cells.Item(8, 10).Value = 1.072820892358421
$ws.Cells.Item(8, 11).Value = 1.069677542046553
$ws.Cells.Item(8, 12).Value = 1.074711861403927
$ws.Cells.Item(8, 13).Value = 1.084131582214221

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.06471586511701
$ws.Cells.Item(9, 4).Value = 1.064410854184381
$ws.Cells.Item(9, 5).Value = 1.069015413720005
$ws.Cells.Item(9, 6).Value = 1.07806366419784
$ws.Cells.Item(9, 9).Value = 1.055339506423568
$ws.Cells.Item(9, 10).Value = 1.070340934330664
$ws.Cells.Item(9, 11).Value = 1.067483454796774
$ws.Cells.Item(9, 12).Value = 1.072073923242188
$ws.Cells.Item(9, 13).Value = 1.081094865979734

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.062521931227502
$ws.Cells.Item(10, 4).Value = 1.062656738991819
$ws.Cells.Item(10, 5).Value = 1.066964169941332
$ws.Cells.Item(10, 6).Value = 1.075750511212564
$ws.Cells.Item(10, 9).Value = 1.054510903166257
$ws.Cells.Item(10, 10).Value = 1.068676265371208
$ws.Cells.Item(10, 11).Value = 1.066009723804426
$ws.Cells.Item(10, 12).Value = 1.070302666590072
$ws.Cells.Item(10, 13).Value = 1.079059849365235

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.061569679915808
$ws.Cells.Item(11, 4).Value = 1.061895216553062
$ws.Cells.Item(11, 5).Value = 1.066073527469757
$ws.Cells.Item(11, 6).Value = 1.074746904804482
$ws.Cells.Item(11, 9).Value = 1.054149375934799
$ws.Cells.Item(11, 10).Value = 1.067952672221851
$ws.Cells.Item(11, 11).Value = 1.065368904136723
$ws.Cells.Item(11, 12).Value = 1.069532610462939
$ws.Cells.Item(11, 13).Value = 1.078176070253452

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.061215622603975
$ws.Cells.Item(12, 4).Value = 1.061612049704285
$ws.Cells.Item(12, 5).Value = 1.065742328453263
$ws.Cells.Item(12, 6).Value = 1.07437381186506
$ws.Cells.Item(12, 9).Value = 1.054014673706316
$ws.Cells.Item(12, 10).Value = 1.067683472977029
$ws.Cells.Item(12, 11).Value = 1.065130465984983
$ws.Cells.Item(12, 12).Value = 1.069246105508211
$ws.Cells.Item(12, 13).Value = 1.077847395950543

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.061291585022746
$ws.Cells.Item(13, 4).Value = 1.061672803788888
$ws.Cells.Item(13, 5).Value = 1.0658133888555
$ws.Cells.Item(13, 6).Value = 1.074453855687239
$ws.Cells.Item(13, 9).Value = 1.054043586628528
$ws.Cells.Item(13, 10).Value = 1.067741236449419
$ws.Cells.Item(13, 11).Value = 1.065181630389944
$ws.Cells.Item(13, 12).Value = 1.069307583266171
$ws.Cells.Item(13, 13).Value = 1.077917915910561

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.061540420593539
$ws.Cells.Item(14, 4).Value = 1.061871816120465
$ws.Cells.Item(14, 5).Value = 1.066046158147571
$ws.Cells.Item(14, 6).Value = 1.074716071177982
$ws.Cells.Item(14, 9).Value = 1.054138249892751
$ws.Cells.Item(14, 10).Value = 1.067930428828082
$ws.Cells.Item(14, 11).Value = 1.065349203146411
$ws.Cells.Item(14, 12).Value = 1.069508937542847
$ws.Cells.Item(14, 13).Value = 1.078148910128409

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.061693689923256
$ws.Cells.Item(15, 4).Value = 1.061994393782385
$ws.Cells.Item(15, 5).Value = 1.066189525089611
$ws.Cells.Item(15, 6).Value = 1.074877589602977
$ws.Cells.Item(15, 9).Value = 1.054196519959848
$ws.Cells.Item(15, 10).Value = 1.068046940038118
$ws.Cells.Item(15, 11).Value = 1.065452395852966
$ws.Cells.Item(15, 12).Value = 1.06963293581071
$ws.Cells.Item(15, 13).Value = 1.078291180118932

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.062585080579933
$ws.Cells.Item(16, 4).Value = 1.062707236506136
$ws.Cells.Item(16, 5).Value = 1.067023226842337
$ws.Cells.Item(16, 6).Value = 1.075817074448948
$ws.Cells.Item(16, 9).Value = 1.054534838559841
$ws.Cells.Item(16, 10).Value = 1.06872422870491
$ws.Cells.Item(16, 11).Value = 1.066052195835897
$ws.Cells.Item(16, 12).Value = 1.070353706927716
$ws.Cells.Item(16, 13).Value = 1.079118447319322

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.063143614263223
$ws.Cells.Item(17, 4).Value = 1.063153849730672
$ws.Cells.Item(17, 5).Value = 1.067545526862393
$ws.Cells.Item(17, 6).Value = 1.076405847799309
$ws.Cells.Item(17, 9).Value = 1.054746321566115
$ws.Cells.Item(17, 10).Value = 1.069148324793548
$ws.Cells.Item(17, 11).Value = 1.066427711678792
$ws.Cells.Item(17, 12).Value = 1.070804994939156
$ws.Cells.Item(17, 13).Value = 1.079636667348151

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.063469179712743
$ws.Cells.Item(18, 4).Value = 1.063414161146432
$ws.Cells.Item(18, 5).Value = 1.067849940717616
$ws.Cells.Item(18, 6).Value = 1.076749077367917
$ws.Cells.Item(18, 9).Value = 1.054869412245329
$ws.Cells.Item(18, 10).Value = 1.069395425050571
$ws.Cells.Item(18, 11).Value = 1.066646485109738
$ws.Cells.Item(18, 12).Value = 1.071067925943708
$ws.Cells.Item(18, 13).Value = 1.079938685685195

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.063580152437737
$ws.Cells.Item(19, 4).Value = 1.063502888482465
$ws.Cells.Item(19, 5).Value = 1.067953698374792
$ws.Cells.Item(19, 6).Value = 1.076866077490465
$ws.Cells.Item(19, 9).Value = 1.054911338358063
$ws.Cells.Item(19, 10).Value = 1.069479634630707
$ws.Cells.Item(19, 11).Value = 1.0667210375091
$ws.Cells.Item(19, 12).Value = 1.071157528360653
$ws.Cells.Item(19, 13).Value = 1.080041623868968

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.063083711475981
$ws.Cells.Item(20, 4).Value = 1.063105952109838
$ws.Cells.Item(20, 5).Value = 1.067489513369828
$ws.Cells.Item(20, 6).Value = 1.076342697911817
$ws.Cells.Item(20, 9).Value = 1.054723658752446
$ws.Cells.Item(20, 10).Value = 1.069102851047567
$ws.Cells.Item(20, 11).Value = 1.066387449174327
$ws.Cells.Item(20, 12).Value = 1.070756606849202
$ws.Cells.Item(20, 13).Value = 1.079581093227884

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.061467154418255
$ws.Cells.Item(21, 4).Value = 1.061813220377338
$ws.Cells.Item(21, 5).Value = 1.065977623802939
$ws.Cells.Item(21, 6).Value = 1.07463886384923
$ws.Cells.Item(21, 9).Value = 1.054110385405138
$ws.Cells.Item(21, 10).Value = 1.067874728162012
$ws.Cells.Item(21, 11).Value = 1.065299868489595
$ws.Cells.Item(21, 12).Value = 1.06944965681242
$ws.Cells.Item(21, 13).Value = 1.078080899180299

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.060448738560575
$ws.Cells.Item(22, 4).Value = 1.060998669520466
$ws.Cells.Item(22, 5).Value = 1.065024866596093
$ws.Cells.Item(22, 6).Value = 1.073565804437943
$ws.Cells.Item(22, 9).Value = 1.053722393629688
$ws.Cells.Item(22, 10).Value = 1.067100098950237
$ws.Cells.Item(22, 11).Value = 1.064613692920313
$ws.Cells.Item(22, 12).Value = 1.068625191857733
$ws.Cells.Item(22, 13).Value = 1.077135352759878

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.060988814690829
$ws.Cells.Item(23, 4).Value = 1.06143064716365
$ws.Cells.Item(23, 5).Value = 1.065530149842573
$ws.Cells.Item(23, 6).Value = 1.074134826342886
$ws.Cells.Item(23, 9).Value = 1.053928304408475
$ws.Cells.Item(23, 10).Value = 1.067510979975346
$ws.Cells.Item(23, 11).Value = 1.06497767424015
$ws.Cells.Item(23, 12).Value = 1.069062517973675
$ws.Cells.Item(23, 13).Value = 1.077636827000119

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.063110779640251
$ws.Cells.Item(24, 4).Value = 1.063127595571906
$ws.Cells.Item(24, 5).Value = 1.067514824182121
$ws.Cells.Item(24, 6).Value = 1.076371233221114
$ws.Cells.Item(24, 9).Value = 1.054733899917606
$ws.Cells.Item(24, 10).Value = 1.069123399503704
$ws.Cells.Item(24, 11).Value = 1.066405642866277
$ws.Cells.Item(24, 12).Value = 1.070778472264133
$ws.Cells.Item(24, 13).Value = 1.079606205555335

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.065564688301232
$ws.Cells.Item(25, 4).Value = 1.065089357342455
$ws.Cells.Item(25, 5).Value = 1.069808741361132
$ws.Cells.Item(25, 6).Value = 1.07895895964056
$ws.Cells.Item(25, 9).Value = 1.055658414171002
$ws.Cells.Item(25, 10).Value = 1.070984040247774
$ws.Cells.Item(25, 11).Value = 1.068052596723486
$ws.Cells.Item(25, 12).Value = 1.072758091355677
$ws.Cells.Item(25, 13).Value = 1.081881756323601
